# Auto-generated Excel COM-interop script that applies the Seraph_Profits market-data refresh.
# For each affected leve row (per-job sheet), updates the price/profit columns H-N
# to the newly-fetched values. Cells with no new value (removed by the refresh) are cleared;
# cells that are newly populated are created.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

$ws.Range("H18").Value = 12499.25
$ws.Range("I18").Value = 19998
$ws.Range("J18").Value = 9999.666999999999
$ws.Range("K18").Value = 19998
$ws.Range("L18").Value = 9999.666999999999
$ws.Range("M18").Value = -19714
$ws.Range("N18").Value = -10567.667

$ws.Range("H21").Value = 11999.5
$ws.Range("I21").Value = 11999.5
$ws.Range("K21").Value = 11999.5
$ws.Range("M21").Value = -11531.5

$ws.Range("H23").Value = 11999.5
$ws.Range("I23").Value = 11999.5
$ws.Range("K23").Value = 11999.5
$ws.Range("M23").Value = -11765.5

$ws.Range("H106").Value = 34511240
$ws.Range("I106").Value = 41694124
$ws.Range("J106").Value = 33399.8
$ws.Range("K106").Value = 41694124
$ws.Range("L106").Value = 33399.8
$ws.Range("M106").Value = -41693493
$ws.Range("N106").Value = -34661.8

$ws.Range("H116").Value = 2500
$ws.Range("I116").Value = 2500
$ws.Range("K116").Value = 2500
$ws.Range("M116").Value = 942

$ws.Range("H137").Value = 3256.8147
$ws.Range("I137").Value = 4057.6155
$ws.Range("J137").Value = 2513.2144
$ws.Range("K137").Value = 12172.8465
$ws.Range("L137").Value = 7539.6432
$ws.Range("M137").Value = -9622.8465
$ws.Range("N137").Value = -12639.6432

$ws.Range("H138").Value = 3104.4695
$ws.Range("I138").Value = 2432.5518
$ws.Range("K138").Value = 7297.655400000001
$ws.Range("M138").Value = -2157.655400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 324319.3
$ws.Range("I32").Value = 1169.2693
$ws.Range("K32").Value = 1169.2693
$ws.Range("M32").Value = -882.2692999999999

$ws.Range("H61").Value = 2517.3333
$ws.Range("I61").Value = 2439.75
$ws.Range("K61").Value = 2439.75
$ws.Range("M61").Value = -2227.75

$ws.Range("H122").Value = 4092.125
$ws.Range("I122").Value = 3237.4285
$ws.Range("K122").Value = 9712.2855
$ws.Range("M122").Value = -7262.2855

$ws.Range("H136").Value = 2517.3333
$ws.Range("I136").Value = 2439.75
$ws.Range("K136").Value = 7319.25
$ws.Range("M136").Value = -4769.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 314
$ws.Range("I5").Value = 314
$ws.Range("K5").Value = 314
$ws.Range("M5").Value = -201

$ws.Range("H82").Value = 185233.17
$ws.Range("J82").Value = 536699.5
$ws.Range("L82").Value = 536699.5
$ws.Range("N82").Value = -537465.5

$ws.Range("H85").Value = 185233.17
$ws.Range("J85").Value = 536699.5
$ws.Range("L85").Value = 536699.5
$ws.Range("N85").Value = -539351.5

$ws.Range("H134").Value = 4264.3335
$ws.Range("I134").Value = 4103
$ws.Range("K134").Value = 12309
$ws.Range("M134").Value = -9774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 9639.333000000001
$ws.Range("J3").Value = 13999
$ws.Range("L3").Value = 13999
$ws.Range("N3").Value = -14225

$ws.Range("H94").Value = 1387.8235
$ws.Range("I94").Value = 726
$ws.Range("K94").Value = 726
$ws.Range("M94").Value = -275

$ws.Range("H134").Value = 4077.0908
$ws.Range("I134").Value = 4205.4443
$ws.Range("K134").Value = 12616.3329
$ws.Range("M134").Value = -10081.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1142.4348
$ws.Range("I11").Value = 339.21738
$ws.Range("K11").Value = 1017.65214
$ws.Range("M11").Value = -877.6521399999999

$ws.Range("H62").Value = 4979.1665
$ws.Range("J62").Value = 4979.1665
$ws.Range("L62").Value = 14937.4995
$ws.Range("N62").Value = -16309.4995

$ws.Range("H65").Value = 4979.1665
$ws.Range("J65").Value = 4979.1665
$ws.Range("L65").Value = 44812.4985
$ws.Range("N65").Value = -51676.4985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 23005.666
$ws.Range("J5").Value = 24507.25
$ws.Range("L5").Value = 24507.25
$ws.Range("N5").Value = -24731.25

$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -26064

$ws.Range("H109").Value = 38349.7
$ws.Range("J109").Value = 38349.7
$ws.Range("L109").Value = 38349.7
$ws.Range("N109").Value = -40429.7

$ws.Range("H132").Value = 2376.5
$ws.Range("I132").Value = 2087.4285
$ws.Range("J132").Value = 4400
$ws.Range("K132").Value = 6262.2855
$ws.Range("L132").Value = 13200
$ws.Range("M132").Value = -3732.2855
$ws.Range("N132").Value = -18260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5005
$ws.Range("J3").Value = 5005
$ws.Range("L3").Value = 5005
$ws.Range("N3").Value = -5229

$ws.Range("H7").Value = 4179
$ws.Range("I7").Value = 4179
$ws.Range("K7").Value = 4179
$ws.Range("M7").Value = -4067

$ws.Range("H15").Value = 5005
$ws.Range("J15").Value = 5005
$ws.Range("L15").Value = 5005
$ws.Range("N15").Value = -5345

$ws.Range("H32").Value = 2266.3333
$ws.Range("I32").Value = 2266.3333
$ws.Range("K32").Value = 2266.3333
$ws.Range("M32").Value = -1949.3333

$ws.Range("H42").Value = 844666.3
$ws.Range("I42").Value = 10999
$ws.Range("K42").Value = 10999
$ws.Range("M42").Value = -10436

$ws.Range("H49").Value = 844666.3
$ws.Range("I49").Value = 10999
$ws.Range("K49").Value = 10999
$ws.Range("M49").Value = -10852

$ws.Range("H122").Value = 7284.5713
$ws.Range("I122").Value = 6298.4
$ws.Range("K122").Value = 18895.2
$ws.Range("M122").Value = -16445.2

$ws.Range("H126").Value = 4179
$ws.Range("I126").Value = 4179
$ws.Range("K126").Value = 12537
$ws.Range("M126").Value = -10067

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 7249.75
$ws.Range("I14").Value = 999.5
$ws.Range("J14").Value = 13500
$ws.Range("K14").Value = 999.5
$ws.Range("L14").Value = 13500
$ws.Range("M14").Value = -831.5
$ws.Range("N14").Value = -13836

$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

$ws.Range("H132").Value = 4831
$ws.Range("I132").Value = 1882.2667
$ws.Range("K132").Value = 5646.800099999999
$ws.Range("M132").Value = -3116.800099999999
